# "Generate Report for Handback" -------------------------------------------
# The localization handback finished: the status moves from "Ready for
# handoff" to "Handed back: in sync with en-US", the per-language "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns get populated (with a link back to the source .md), and the
# columns that now hold long file names/links are widened.

$wb = $excel.ActiveWorkbook

$newStatus  = "Handed back: in sync with en-US"
$mdName     = "ae134d18-3dda-42d6-a21a-1ee067d5c0e3.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7eb0d0c03cbfcedd5d850c8a9c9dadef573e9b34/e2e/ae134d18-3dda-42d6-a21a-1ee067d5c0e3.md"
$zhXlf      = "ae134d18-3dda-42d6-a21a-1ee067d5c0e3.f11f31b3b97422a3ea8d9b73dc1178f22f74256d.zh-cn.xlf"
$deXlf      = "ae134d18-3dda-42d6-a21a-1ee067d5c0e3.f11f31b3b97422a3ea8d9b73dc1178f22f74256d.de-de.xlf"
$zhHandback = "2016-08-24 04:56:28"
$deHandback = "2016-08-24 04:56:35"

# A ColumnWidth value that round-trips (through the host's character-grid
# rounding) to the widened "40"-style columns used for long file names.
$wideWidth = 39.16666666666667
# ... and one that round-trips close to the ~30-character width now used
# for the Status / zh-cn / de-de columns.
$midWidth  = 29.16666666666666

# ---------------------------------------------------------------------
# Overview sheet: status summary + widen the zh-cn/de-de columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $midWidth
$wsOverview.Columns.Item(6).ColumnWidth = $midWidth

# ---------------------------------------------------------------------
# Helper to apply the per-language handback report to a language sheet
# ---------------------------------------------------------------------
function Set-HandbackReport($wsName, $xlfName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($wsName)

    # Status -> handed back
    $ws.Range("C2").Value = $newStatus

    # Latest Target File: the source markdown file, linked back to GitHub
    $ws.Range("I2").Value = $mdName
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl, "", "", $mdName) | Out-Null
    $ws.Range("I2").Font.Underline = $true
    $ws.Range("I2").Font.Color = 15570276

    # Latest Handback File: the generated xliff for this language
    $ws.Range("J2").Value = $xlfName

    # Latest Handback DateTime
    $ws.Range("K2").Value = $handbackDateTime

    # Widen the columns that now carry long file names / links
    $ws.Columns.Item(3).ColumnWidth = $midWidth
    $ws.Columns.Item(9).ColumnWidth = $wideWidth
    $ws.Columns.Item(10).ColumnWidth = $wideWidth
}

Set-HandbackReport "zh-cn" $zhXlf $zhHandback
Set-HandbackReport "de-de" $deXlf $deHandback
